$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 21-25 (IdleBuilder, IdleConstruct, IdleControl, IdleRepair, IdleSummonEntity)
# now have a deconstructor implemented: "Has Deconstruct?" -> Yes,
# "Lifetime" -> Self Limiting, "Deconstructed At" -> Self.
$rows = 21..25
foreach ($r in $rows) {
    $ws.Range("D$r").Value = "Yes"
    $ws.Range("E$r").Value = "Self Limiting"
    $ws.Range("F$r").Value = "Self"
}

# Match the highlighting style used for other "Yes" rows (e.g. row 2) by
# copying its formatting (fill colour) onto the updated rows.
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A21:G25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Move the active selection to J7.
$ws.Range("J7").Select() | Out-Null
